# Commit: "shorten names of tabs"
#
# 1. Move the "birth outcome distribution" sheet tab so it sits right after
#    "distributions" (i.e. becomes the 4th tab, right before "RRStunting").
# 2. Shorten two long tab names:
#      "OR stunting for complementary f" -> "OR stunting for complements"
#      "OR exclusive breastfeeding by p" -> "OR exclusiveBF by intervention"

$wb = $excel.ActiveWorkbook

# --- 1. Reorder: move "birth outcome distribution" before "RRStunting" ---
$moving = $wb.Worksheets.Item("birth outcome distribution")
$anchor = $wb.Worksheets.Item("RRStunting")
$moving.Move($anchor)

# --- 2. Rename tabs ---
$wb.Worksheets.Item("OR stunting for complementary f").Name = "OR stunting for complements"
$wb.Worksheets.Item("OR exclusive breastfeeding by p").Name = "OR exclusiveBF by intervention"
